# Applies the "trade closed / new trades opened" update to the live
# trading results workbook.
#
# Summary of changes:
#  - Summary sheet: refresh aggregate metrics after closing trade #32
#    (MarketMaking) and opening two new trades.
#  - Strategy Status sheet: refresh MarketMaking strategy row.
#  - All Trades sheet: close trade #32 (row 33) and append two new
#    OPEN trades (rows 62 and 63).
#  - momentum sheet: append the new OPEN momentum trade.
#  - HighProbConvergence sheet: append the new OPEN HighProbConvergence trade.
#  - MarketMaking sheet: close its trade (row 4, mirrors All Trades row 33).

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($Cell, [string]$Text)
    # Force the cell to be treated as literal text so values like
    # "2026-02-18" or "00:09:27" are not auto-converted into Excel
    # date/time serial numbers.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$sumWs = $wb.Worksheets.Item("Summary")
$sumWs.Cells.Item(3, 2).Value = 1499.88   # Current Capital
$sumWs.Cells.Item(4, 2).Value = 0.98      # Total P&L $
$sumWs.Cells.Item(5, 2).Value = 0.61      # Total P&L %
$sumWs.Cells.Item(6, 2).Value = 32        # Total Trades
$sumWs.Cells.Item(7, 2).Value = 17        # Winning Trades
$sumWs.Cells.Item(9, 2).Value = 53.12     # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet (row 6 = MarketMaking)
# ---------------------------------------------------------------
$statusWs = $wb.Worksheets.Item("Strategy Status")
$statusWs.Cells.Item(6, 3).Value = 99.88                 # Capital
$statusWs.Cells.Item(6, 4).Value = 3                     # Trades
$statusWs.Cells.Item(6, 5).Value = 0.07000000000000001   # P&L $
$statusWs.Cells.Item(6, 6).Value = -0.12                 # P&L %
$statusWs.Cells.Item(6, 7).Value = 33.33                 # Win Rate %

# ---------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------
$allWs = $wb.Worksheets.Item("All Trades")

# Close trade #32 (MarketMaking) on row 33.
$allWs.Cells.Item(33, 7).Value = 0.727221     # Exit Price
$allWs.Cells.Item(33, 8).Value = "CLOSED"     # Status
$allWs.Cells.Item(33, 9).Value = 11.8801      # P&L %
$allWs.Cells.Item(33, 10).Value = 0.08        # P&L $
$allWs.Cells.Item(33, 11).Value = 99.88       # Capital After
$allWs.Cells.Item(33, 12).Value = "early_exit" # Exit Reason
$allWs.Cells.Item(33, 13).Value = 0.19        # Duration (min)

# New trade #61 -> row 62 (momentum, OPEN)
Set-TextCell $allWs.Cells.Item(62, 2) "2026-02-18"
Set-TextCell $allWs.Cells.Item(62, 3) "00:09:27"
$allWs.Cells.Item(62, 1).Value = 61
$allWs.Cells.Item(62, 4).Value = "momentum"
$allWs.Cells.Item(62, 5).Value = "DOWN"
$allWs.Cells.Item(62, 6).Value = 0.65
$allWs.Cells.Item(62, 8).Value = "OPEN"
$allWs.Cells.Item(62, 9).Value = 0
$allWs.Cells.Item(62, 10).Value = 0
$allWs.Cells.Item(62, 11).Value = 100
$allWs.Cells.Item(62, 13).Value = 0
$allWs.Cells.Item(62, 14).Value = 0
$allWs.Cells.Item(62, 15).Value = 0
$allWs.Cells.Item(62, 16).Value = 0.9
$allWs.Cells.Item(62, 17).Value = "Downward momentum: -17.822% over 10 samples"

# New trade #62 -> row 63 (HighProbConvergence, OPEN)
Set-TextCell $allWs.Cells.Item(63, 2) "2026-02-18"
Set-TextCell $allWs.Cells.Item(63, 3) "00:09:28"
$allWs.Cells.Item(63, 1).Value = 62
$allWs.Cells.Item(63, 4).Value = "HighProbConvergence"
$allWs.Cells.Item(63, 5).Value = "UP"
$allWs.Cells.Item(63, 6).Value = 0.35
$allWs.Cells.Item(63, 8).Value = "OPEN"
$allWs.Cells.Item(63, 9).Value = 0
$allWs.Cells.Item(63, 10).Value = 0
$allWs.Cells.Item(63, 11).Value = 100
$allWs.Cells.Item(63, 13).Value = 0
$allWs.Cells.Item(63, 14).Value = 0
$allWs.Cells.Item(63, 15).Value = 0
$allWs.Cells.Item(63, 16).Value = 0.95
$allWs.Cells.Item(63, 17).Value = "Mean reversion UP: price 16.33% below mean (z=-3.00)"

# ---------------------------------------------------------------
# momentum sheet: append new OPEN trade #61 as row 3
# ---------------------------------------------------------------
$momWs = $wb.Worksheets.Item("momentum")
Set-TextCell $momWs.Cells.Item(3, 2) "2026-02-18"
Set-TextCell $momWs.Cells.Item(3, 3) "00:09:27"
$momWs.Cells.Item(3, 1).Value = 61
$momWs.Cells.Item(3, 4).Value = "momentum"
$momWs.Cells.Item(3, 5).Value = "DOWN"
$momWs.Cells.Item(3, 6).Value = 0.65
$momWs.Cells.Item(3, 8).Value = "OPEN"
$momWs.Cells.Item(3, 9).Value = 0
$momWs.Cells.Item(3, 10).Value = 0
$momWs.Cells.Item(3, 11).Value = 100
$momWs.Cells.Item(3, 12).Value = 0
$momWs.Cells.Item(3, 13).Value = 0
$momWs.Cells.Item(3, 14).Value = 0.9
$momWs.Cells.Item(3, 15).Value = "Downward momentum: -17.822% over 10 samples"
$momWs.Cells.Item(3, 17).Value = 0

# ---------------------------------------------------------------
# HighProbConvergence sheet: append new OPEN trade #62 as row 3
# ---------------------------------------------------------------
$hpcWs = $wb.Worksheets.Item("HighProbConvergence")
Set-TextCell $hpcWs.Cells.Item(3, 2) "2026-02-18"
Set-TextCell $hpcWs.Cells.Item(3, 3) "00:09:28"
$hpcWs.Cells.Item(3, 1).Value = 62
$hpcWs.Cells.Item(3, 4).Value = "HighProbConvergence"
$hpcWs.Cells.Item(3, 5).Value = "UP"
$hpcWs.Cells.Item(3, 6).Value = 0.35
$hpcWs.Cells.Item(3, 8).Value = "OPEN"
$hpcWs.Cells.Item(3, 9).Value = 0
$hpcWs.Cells.Item(3, 10).Value = 0
$hpcWs.Cells.Item(3, 11).Value = 100
$hpcWs.Cells.Item(3, 12).Value = 0
$hpcWs.Cells.Item(3, 13).Value = 0
$hpcWs.Cells.Item(3, 14).Value = 0.95
$hpcWs.Cells.Item(3, 15).Value = "Mean reversion UP: price 16.33% below mean (z=-3.00)"
$hpcWs.Cells.Item(3, 17).Value = 0

# ---------------------------------------------------------------
# MarketMaking sheet: close trade on row 4 (mirrors All Trades row 33)
# ---------------------------------------------------------------
$mmWs = $wb.Worksheets.Item("MarketMaking")
$mmWs.Cells.Item(4, 7).Value = 0.727221       # Exit Price
$mmWs.Cells.Item(4, 8).Value = "CLOSED"       # Status
$mmWs.Cells.Item(4, 9).Value = 11.8801        # P&L %
$mmWs.Cells.Item(4, 10).Value = 0.08          # P&L $
$mmWs.Cells.Item(4, 11).Value = 99.88         # Capital After
$mmWs.Cells.Item(4, 16).Value = "early_exit"  # Exit Reason
$mmWs.Cells.Item(4, 17).Value = 0.19          # Duration (min)
